$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the student names in column B (NAME) with single-letter placeholders,
# keeping the COURSE column (C) values unchanged for each row.
$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "B"
$ws.Range("B4").Value = "C"
$ws.Range("B5").Value = "D"
